$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: id, MovieId, Title, Rating, Review, Year
# (a new "id" (GUID) column is inserted before the old "ID" column, which is
# renamed "MovieId"; Title/Rating/Review/Year all shift one column to the right)

$data = @(
    @("id",                                   "MovieId", "Title",                                               "Rating", "Review",                                      "Year"),
    @("d861be4a-de63-49ba-94e0-57486b060d90",  1,         "The Lord of the Rings: The Fellowship of the Ring",  "PG-13",  "Long, but good",                              2001),
    @("bd7d27f8-2f3d-4044-8117-e9e71e351339",  2,         "The Lord of the Rings: The Two Towers",              "PG-13",  "Needs resolution",                            2002),
    @("a7730d44-e048-4879-b571-a8a92a94c1be",  3,         "The Lord of the Rings: The Return of the King",      "PG-13",  "Really Good",                                 2003),
    @("de44443e-7c36-4a51-8101-be42d0b572a1",  4,         "Top Gun",                                             "PG",     "A favorite",                                  1986),
    @("86ae4a99-30aa-42f3-bf6e-0d08e535ff7d",  5,         "Top Gun: Maverick",                                   "PG-13",  "One of the best movies ever",                 2022),
    @("a264bd90-22d5-47b7-aefe-0f6df48de7a3",  6,         "Apollo 13",                                           "PG-13",  "We have a problem, but it's not the movie",  1995),
    @("572eedce-7e46-4d32-915a-f07c529fed2d",  7,         "Bill and Ted's Excellent Adventure",                  "PG",     "Excellent",                                   1989)
)

# Clear any leftover cell styling from the old layout (old A1:E8 cells used
# style index 1) before writing the new values - the new layout uses the
# default style everywhere.
$ws.Range("A1:F8").ClearFormats()

# Match the original authoring order: the "MovieId" header (renamed from the
# old "ID" header) was typed first, then the new "id" header, then the new
# GUID column was filled in top-to-bottom, and finally the rest of the
# (unchanged) table content was rewritten in normal row-major order.
$ws.Cells.Item(1, 2).Value = $data[0][1]
$ws.Cells.Item(1, 1).Value = $data[0][0]
for ($r = 1; $r -lt $data.Length; $r++) {
    $ws.Cells.Item($r + 1, 1).Value = $data[$r][0]
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        if ($r -eq 0 -and $c -le 1) { continue }
        if ($r -ge 1 -and $c -eq 0) { continue }
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Old column E (unused now that the old E header/data moved to F) no longer
# carries any value beyond row 8, nothing further to clear there.

# Column widths: new id column gets its own bestfit-ish width, Title/Rating/
# Review shift from B/C/D to C/D/E keeping their original widths, column B
# (MovieId) keeps the default width.
$ws.Columns("A").ColumnWidth = 38.140625
$ws.Columns("C").ColumnWidth = 45
$ws.Columns("D").ColumnWidth = 13.42578125
$ws.Columns("E").ColumnWidth = 26.5703125

$ws.Range("C9").Select()

Write-Output "done"
